# "Nuevo registro" - opcion 2 funcional, se registra un nuevo libro
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a couple of existing author names (C2, C6)
$ws.Range("C6").Value = "juancito"
$ws.Range("C2").Value = "asdf"

# Register new book records in the next free rows
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "libro"
$ws.Range("C16").Value = "autor"
$ws.Range("D16").Value = 20

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "azul"
$ws.Range("C17").Value = "verde"
$ws.Range("D17").Value = 500
